$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.033963283021942
$ws.Range("D2").Value = 1.036364150783965
$ws.Range("E2").Value = 1.042777961901193
$ws.Range("F2").Value = 1.052831699119817
$ws.Range("I2").Value = 1.033085759734932
$ws.Range("J2").Value = 1.03908512900133
$ws.Range("K2").Value = 1.03915810398726
$ws.Range("L2").Value = 1.045553689824181
$ws.Range("M2").Value = 1.055579343117329
$ws.Range("N2").Value = 1.016918472441336

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.035269566251127
$ws.Range("D3").Value = 1.037350660930939
$ws.Range("E3").Value = 1.044011395080187
$ws.Range("F3").Value = 1.054320994020457
$ws.Range("I3").Value = 1.033374583315384
$ws.Range("J3").Value = 1.04003253949056
$ws.Range("K3").Value = 1.039953730936244
$ws.Range("L3").Value = 1.046596917790977
$ws.Range("M3").Value = 1.056879828409117
$ws.Range("N3").Value = 1.01723796926151

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.036113699194543
$ws.Range("D4").Value = 1.037987564845148
$ws.Range("E4").Value = 1.04480888915165
$ws.Range("F4").Value = 1.055284439131023
$ws.Range("I4").Value = 1.033559086338679
$ws.Range("J4").Value = 1.040644011864765
$ws.Range("K4").Value = 1.040466513564812
$ws.Range("L4").Value = 1.047270779300888
$ws.Range("M4").Value = 1.057720594144355
$ws.Range("N4").Value = 1.017444027952531

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.036468308920704
$ws.Range("D5").Value = 1.038254978207963
$ws.Range("E5").Value = 1.045144010821499
$ws.Range("F5").Value = 1.05568942151139
$ws.Range("I5").Value = 1.033636081271939
$ws.Range("J5").Value = 1.040900702898949
$ws.Range("K5").Value = 1.040681600887011
$ws.Range("L5").Value = 1.047553792150664
$ws.Range("M5").Value = 1.058073880170547
$ws.Range("N5").Value = 1.017530493817832

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.036527834110921
$ws.Range("D6").Value = 1.038299858155343
$ws.Range("E6").Value = 1.045200270832815
$ws.Range("F6").Value = 1.055757417042656
$ws.Range("I6").Value = 1.033648975651888
$ws.Range("J6").Value = 1.040943780723966
$ws.Range("K6").Value = 1.040717686529286
$ws.Range("L6").Value = 1.047601295012581
$ws.Range("M6").Value = 1.058133188488609
$ws.Range("N6").Value = 1.017545002389612

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.036118438538176
$ws.Range("D7").Value = 1.037991139372862
$ws.Range("E7").Value = 1.04481336762858
$ws.Range("F7").Value = 1.055289850715828
$ws.Range("I7").Value = 1.033560117387952
$ws.Range("J7").Value = 1.040647443240946
$ws.Range("K7").Value = 1.04046938948108
$ws.Range("L7").Value = 1.04727456202294
$ws.Range("M7").Value = 1.057725315439838
$ws.Range("N7").Value = 1.017445183945993

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.034404981973263
$ws.Range("D8").Value = 1.036697844835253
$ws.Range("E8").Value = 1.043194936847198
$ws.Range("F8").Value = 1.053335063409497
$ws.Range("I8").Value = 1.033183863405718
$ws.Range("J8").Value = 1.039405636692713
$ws.Range("K8").Value = 1.039427413802104
$ws.Range("L8").Value = 1.045906499176015
$ws.Range("M8").Value = 1.056019002927847
$ws.Range("N8").Value = 1.017026588629381

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.031376884466236
$ws.Range("D9").Value = 1.034407807086443
$ws.Range("E9").Value = 1.040338150781089
$ws.Range("F9").Value = 1.049888518530873
$ws.Range("I9").Value = 1.032502548460414
$ws.Range("J9").Value = 1.037205295935338
$ws.Range("K9").Value = 1.037575587565116
$ws.Range("L9").Value = 1.043486641593236
$ws.Range("M9").Value = 1.053006445988436
$ws.Range("N9").Value = 1.01628374387137

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.029351985096348
$ws.Range("D10").Value = 1.032873511607486
$ws.Range("E10").Value = 1.038430105303647
$ws.Range("F10").Value = 1.047589205950183
$ws.Range("I10").Value = 1.032035971409568
$ws.Range("J10").Value = 1.035730071293444
$ws.Range("K10").Value = 1.036330318074237
$ws.Range("L10").Value = 1.041867043807881
$ws.Range("M10").Value = 1.050993902329774
$ws.Range("N10").Value = 1.015784942566859

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.02847365693921
$ws.Range("D11").Value = 1.032207306956935
$ws.Range("E11").Value = 1.037603014611959
$ws.Range("F11").Value = 1.046593127816172
$ws.Range("I11").Value = 1.031830991399079
$ws.Range("J11").Value = 1.035089265254123
$ws.Range("K11").Value = 1.035788527838491
$ws.Range("L11").Value = 1.041164185680429
$ws.Range("M11").Value = 1.050121400058778
$ws.Range("N11").Value = 1.015568095856215

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028147170981265
$ws.Range("D12").Value = 1.031959568676464
$ws.Range("E12").Value = 1.037295657461782
$ws.Range("F12").Value = 1.046223064974262
$ws.Range("I12").Value = 1.031754408356949
$ws.Range("J12").Value = 1.03485093358484
$ws.Range("K12").Value = 1.035586892235397
$ws.Range("L12").Value = 1.040902874199789
$ws.Range("M12").Value = 1.049797149797999
$ws.Range("N12").Value = 1.015487418585528

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028217214083889
$ws.Range("D13").Value = 1.032012722108564
$ws.Range("E13").Value = 1.037361592968761
$ws.Range("F13").Value = 1.046302448201921
$ws.Range("I13").Value = 1.031770855807422
$ws.Range("J13").Value = 1.034902070523894
$ws.Range("K13").Value = 1.035630161478817
$ws.Range("L13").Value = 1.040958937265674
$ws.Range("M13").Value = 1.049866710087634
$ws.Range("N13").Value = 1.015504730069488

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028446674349158
$ws.Range("D14").Value = 1.032186834567553
$ws.Range("E14").Value = 1.037577611232315
$ws.Range("F14").Value = 1.046562539871771
$ws.Range("I14").Value = 1.031824670092761
$ws.Range("J14").Value = 1.035069570983932
$ws.Range("K14").Value = 1.035771868561493
$ws.Range("L14").Value = 1.041142590484493
$ws.Range("M14").Value = 1.050094600808117
$ws.Range("N14").Value = 1.01556142972188

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028588020987292
$ws.Range("D15").Value = 1.032294073776558
$ws.Range("E15").Value = 1.037710688682026
$ws.Range("F15").Value = 1.046722780809557
$ws.Range("I15").Value = 1.03185776793229
$ws.Range("J15").Value = 1.035172732651671
$ws.Range("K15").Value = 1.035859127127648
$ws.Range("L15").Value = 1.041255713544736
$ws.Range("M15").Value = 1.050234989906063
$ws.Range("N15").Value = 1.015596346888494

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029410244076966
$ws.Range("D16").Value = 1.032917686290455
$ws.Range("E16").Value = 1.038484977364274
$ws.Range("F16").Value = 1.047655302170626
$ws.Range("I16").Value = 1.032049513015865
$ws.Range("J16").Value = 1.035772556521843
$ws.Range("K16").Value = 1.036366220294669
$ws.Range("L16").Value = 1.041913656939655
$ws.Range("M16").Value = 1.051051784708627
$ws.Range("N16").Value = 1.015799315698042

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029925588570453
$ws.Range("D17").Value = 1.033308365931413
$ws.Range("E17").Value = 1.038970425596472
$ws.Range("F17").Value = 1.048240120731226
$ws.Range("I17").Value = 1.032168999196861
$ws.Range("J17").Value = 1.036148265704165
$ws.Range("K17").Value = 1.036683613561138
$ws.Range("L17").Value = 1.042325946503937
$ws.Range("M17").Value = 1.051563851717493
$ws.Range("N17").Value = 1.015926400941103

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030226032714666
$ws.Range("D18").Value = 1.033536064991793
$ws.Range("E18").Value = 1.039253493344092
$ws.Range("F18").Value = 1.048581191460245
$ws.Range("I18").Value = 1.032238409010066
$ws.Range("J18").Value = 1.036367215475283
$ws.Range("K18").Value = 1.036868494919225
$ws.Range("L18").Value = 1.042566277636813
$ws.Range("M18").Value = 1.051862430051129
$ws.Range("N18").Value = 1.016000444479054

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030328451533984
$ws.Range("D19").Value = 1.033613674427493
$ws.Range("E19").Value = 1.039349997679198
$ws.Range("F19").Value = 1.048697480470935
$ws.Range("I19").Value = 1.032262027752288
$ws.Range("J19").Value = 1.036441838717836
$ws.Range("K19").Value = 1.036931492577335
$ws.Range("L19").Value = 1.042648198955025
$ws.Range("M19").Value = 1.051964220454451
$ws.Range("N19").Value = 1.016025677346241

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029870312288961
$ws.Range("D20").Value = 1.033266468113732
$ws.Range("E20").Value = 1.038918350523733
$ws.Range("F20").Value = 1.048177379834024
$ws.Range("I20").Value = 1.032156208892676
$ws.Range("J20").Value = 1.036107975866148
$ws.Range("K20").Value = 1.03664958602045
$ws.Range("L20").Value = 1.042281727322217
$ws.Range("M20").Value = 1.05150892233073
$ws.Range("N20").Value = 1.015912774505094

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.02837911057224
$ws.Range("D21").Value = 1.032135570548127
$ws.Range("E21").Value = 1.037514003106539
$ws.Range("F21").Value = 1.046485951477322
$ws.Range("I21").Value = 1.031808835403293
$ws.Range("J21").Value = 1.035020254799129
$ws.Range("K21").Value = 1.035730150163311
$ws.Range("L21").Value = 1.041088515810447
$ws.Range("M21").Value = 1.05002749724184
$ws.Range("N21").Value = 1.015544736706189

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027440166055565
$ws.Range("D22").Value = 1.031422907600754
$ws.Range("E22").Value = 1.036630228619624
$ws.Range("F22").Value = 1.045422046055396
$ws.Range("I22").Value = 1.031587856447116
$ws.Range("J22").Value = 1.034334578856994
$ws.Range("K22").Value = 1.035149803066335
$ws.Range("L22").Value = 1.040336913425927
$ws.Range("M22").Value = 1.049095114647217
$ws.Range("N22").Value = 1.015312579683168

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027938049680012
$ws.Range("D23").Value = 1.031800858569786
$ws.Range("E23").Value = 1.037098812050481
$ws.Range("F23").Value = 1.045986085867717
$ws.Range("I23").Value = 1.031705245747469
$ws.Range("J23").Value = 1.034698238849384
$ws.Range("K23").Value = 1.035457671369204
$ws.Range("L23").Value = 1.040735484522784
$ws.Range("M23").Value = 1.049589480263157
$ws.Range("N23").Value = 1.015435722667444

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029895289716894
$ws.Range("D24").Value = 1.033285400483212
$ws.Range("E24").Value = 1.038941881273023
$ws.Range("F24").Value = 1.048205729880548
$ws.Range("I24").Value = 1.032161989159451
$ws.Range("J24").Value = 1.036126181712125
$ws.Range("K24").Value = 1.036664962369888
$ws.Range("L24").Value = 1.042301708531507
$ws.Range("M24").Value = 1.051533742870091
$ws.Range("N24").Value = 1.015918931961799

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.032160787256604
$ws.Range("D25").Value = 1.03500116589195
$ws.Range("E25").Value = 1.04107730260154
$ws.Range("F25").Value = 1.050779795845718
$ws.Range("I25").Value = 1.032680860151608
$ws.Range("J25").Value = 1.037775591628575
$ws.Range("K25").Value = 1.038056207802188
$ws.Range("L25").Value = 1.044113338662731
$ws.Range("M25").Value = 1.053785980274965
$ws.Range("N25").Value = 1.016476412289388
